$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 84
$ws.Range("A84").Value = 8
$ws.Range("B84").Value = "Terminal La Palmera de La Serena"
$ws.Range("C84").Value = "Coquimbo"
$ws.Range("D84").Value = 44911
$ws.Range("D84").NumberFormat = $ws.Range("D83").NumberFormat
$ws.Range("E84").Value = 4
$ws.Range("F84").Value = "Fruta"
$ws.Range("G84").Value = 100103
$ws.Range("H84").Value = "Frutos de hueso (carozo)"
$ws.Range("I84").Value = 100103003
$ws.Range("J84").Value = "Damasco"
$ws.Range("K84").Value = "Dina"
$ws.Range("L84").Value = "Especial"
$ws.Range("M84").Value = 200
$ws.Range("N84").Value = 22000
$ws.Range("O84").Value = 23000
$ws.Range("P84").Value = 22500
$ws.Range("Q84").Value = "$/caja 16 kilos"
$ws.Range("R84").Value = "Región Metropolitana"
$ws.Range("S84").Value = 1406
$ws.Range("T84").Value = 16

# Row 85
$ws.Range("A85").Value = 8
$ws.Range("B85").Value = "Terminal La Palmera de La Serena"
$ws.Range("C85").Value = "Coquimbo"
$ws.Range("D85").Value = 44911
$ws.Range("D85").NumberFormat = $ws.Range("D83").NumberFormat
$ws.Range("E85").Value = 4
$ws.Range("F85").Value = "Fruta"
$ws.Range("G85").Value = 100103
$ws.Range("H85").Value = "Frutos de hueso (carozo)"
$ws.Range("I85").Value = 100103003
$ws.Range("J85").Value = "Damasco"
$ws.Range("K85").Value = "Dina"
$ws.Range("L85").Value = "Primera"
$ws.Range("M85").Value = 140
$ws.Range("N85").Value = 19000
$ws.Range("O85").Value = 20000
$ws.Range("P85").Value = 19500
$ws.Range("Q85").Value = "$/caja 16 kilos"
$ws.Range("R85").Value = "Región Metropolitana"
$ws.Range("S85").Value = 1219
$ws.Range("T85").Value = 16
